# Regen save_data to use K (strikeouts) instead of Strike# in column G.
# Updates the already-computed K values for each game row (2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 4
    4  = 6
    5  = 13
    6  = 9
    7  = 6
    8  = 7
    9  = 7
    10 = 9
    11 = 4
    12 = 10
    13 = 4
    14 = 3
    15 = 2
    16 = 1
    17 = 3
    18 = 4
    19 = 2
    20 = 3
    21 = 2
    22 = 2
    23 = 4
    24 = 3
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
